$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# Headers: both "default" (Item 1) and "first page" (Item 2) headers contain
# the BTec logo picture, whose docPr/cNvPr "name" label needs to flip from
# "image2.jpg" to "image1.jpg".
for ($i = 1; $i -le 2; $i++) {
    $hdr = $sec.Headers.Item($i)
    if ($hdr.Exists) {
        $shapes = $hdr.Range.InlineShapes
        for ($j = 1; $j -le $shapes.Count; $j++) {
            $shp = $shapes.Item($j)
            if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
                $shp.Name = "image1.jpg"
            }
        }
    }
}

# Footers: both "default" (Item 1) and "first page" (Item 2) footers contain
# the Pearson logo picture, whose docPr/cNvPr "name" label needs to flip from
# "image1.png" to "image2.png".
for ($i = 1; $i -le 2; $i++) {
    $ftr = $sec.Footers.Item($i)
    if ($ftr.Exists) {
        $shapes = $ftr.Range.InlineShapes
        for ($j = 1; $j -le $shapes.Count; $j++) {
            $shp = $shapes.Item($j)
            if ($shp.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
                $shp.Name = "image2.png"
            }
        }
    }
}
